$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume data range to Text format so numeric-looking
# strings (e.g. "1.000", "234.91") are written verbatim and not coerced
# into floating-point numbers. Style is reset to Normal afterwards so no
# visible formatting change remains.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.652.48'
$ws.Range('E2').Value = '  +1.39%  '

$ws.Range('D3').Value = '1.861.64'
$ws.Range('E3').Value = '  +0.27%  '

$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '234.91'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('D7').Value = '0.4715'
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '0.2752'
$ws.Range('E8').Value = '  +0.67%  '

$ws.Range('D9').Value = '0.06356'
$ws.Range('E9').Value = '  -0.96%  '

$ws.Range('D10').Value = '17.58'
$ws.Range('E10').Value = '  +8.60%  '

$ws.Range('D11').Value = '1.855.75'
$ws.Range('E11').Value = '  +0.00%  '

$ws.Range('D12').Value = '0.07441'
$ws.Range('E12').Value = '  -0.06%  '

$ws.Range('D13').Value = '5.221'
$ws.Range('E13').Value = '  +4.57%  '

$ws.Range('D14').Value = '85.04'
$ws.Range('E14').Value = '  -0.34%  '

$ws.Range('D15').Value = '0.6319'
$ws.Range('E15').Value = '  +0.21%  '

$ws.Range('D16').Value = '30.629.14'
$ws.Range('E16').Value = '  +1.45%  '

$ws.Range('D17').Value = '243.52'
$ws.Range('E17').Value = '  +4.87%  '

$ws.Range('D19').Value = '12.83'
$ws.Range('E19').Value = '  +0.88%  '

$ws.Range('D20').Value = '0.000007369'
$ws.Range('E20').Value = '  +0.15%  '

$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.07%  '

$ws.Range('D22').Value = '4.997'
$ws.Range('E22').Value = '  -0.61%  '

$ws.Range('D23').Value = '6.044'

$ws.Range('D24').Value = '9.342'
$ws.Range('E24').Value = '  +1.09%  '

$ws.Range('D25').Value = '165.01'
$ws.Range('E25').Value = '  +0.10%  '

$ws.Range('D26').Value = '18.12'
$ws.Range('E26').Value = '  +1.46%  '

$ws.Range('D27').Value = '1.890'
$ws.Range('E27').Value = '  +0.31%  '

$ws.Range('E28').Value = '  -0.73%  '

$ws.Range('D29').Value = '1.382'
$ws.Range('E29').Value = '  +0.10%  '

$ws.Range('D30').Value = '4.068'
$ws.Range('E30').Value = '  -1.52%  '

$ws.Range('D31').Value = '3.868'
$ws.Range('E31').Value = '  -1.41%  '

$ws.Range('D32').Value = '0.04928'
$ws.Range('E32').Value = '  +0.70%  '

$ws.Range('E33').Value = '  +0.68%  '

$ws.Range('D34').Value = '0.7065'
$ws.Range('E34').Value = '  -2.16%  '

$ws.Range('D35').Value = '2.711'
$ws.Range('E35').Value = '  +0.70%  '

$ws.Range('D36').Value = '0.01918'
$ws.Range('E36').Value = '  +1.07%  '

$ws.Range('D37').Value = '2.688'
$ws.Range('E37').Value = '  +1.89%  '

$ws.Range('D38').Value = '0.8803'
$ws.Range('E38').Value = '  -2.69%  '

$ws.Range('D39').Value = '1.999'
$ws.Range('E39').Value = '  +1.44%  '

$ws.Range('E40').Value = '  -0.23%  '

$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  +0.18%  '

$ws.Range('D42').Value = '5.548'
$ws.Range('E42').Value = '  +0.63%  '

$ws.Range('D43').Value = '0.4081'
$ws.Range('E43').Value = '  -0.49%  '

$ws.Range('D44').Value = '7.277'
$ws.Range('E44').Value = '  +2.32%  '

$ws.Range('D45').Value = '63.29'
$ws.Range('E45').Value = '  +3.75%  '

$ws.Range('D46').Value = '0.1216'
$ws.Range('E46').Value = '  +1.57%  '

$ws.Range('E47').Value = '  +1.06%  '

$ws.Range('D48').Value = '8.598'
$ws.Range('E48').Value = '  -1.11%  '

$ws.Range('D49').Value = '0.05545'
$ws.Range('E49').Value = '  -0.44%  '

$ws.Range('E50').Value = '  -2.37%  '

$ws.Range('D51').Value = '0.3698'
$ws.Range('E51').Value = '  +0.11%  '

$dataRange.Style = "Normal"
